$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 128, pushing existing rows 128:257 down to 129:258
$ws.Rows(128).Insert()

# Populate the new row 128 with the new data point
$ws.Range("A128").Value = 10
$ws.Range("B128").Value = "Vega Modelo de Temuco"
$ws.Range("C128").Value = "La Araucanía"
$ws.Range("D128").Value = 44781
$ws.Range("E128").Value = 9
$ws.Range("F128").Value = 100112052
$ws.Range("G128").Value = "Albahaca"
$ws.Range("H128").Value = "Sin especificar"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 80
$ws.Range("K128").Value = 6000
$ws.Range("L128").Value = 6000
$ws.Range("M128").Value = 6000
$ws.Range("N128").Value = "$/paquete"
$ws.Range("O128").Value = "Región de Arica y Parinacota"
$ws.Range("P128").Value = 6000
$ws.Range("Q128").Value = 1
$ws.Range("R128").Value = "Hortaliza"
